$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")

# Row 2
$ws.Range("H2").Value = 1025.258
$ws.Range("I2").Value = 219.88235
$ws.Range("J2").Value = 2003.2142
$ws.Range("K2").Value = 219.88235
$ws.Range("L2").Value = 2003.2142
$ws.Range("M2").Value = -106.88235
$ws.Range("N2").Value = -2229.2142

# Row 20
$ws.Range("H20").Value = 1700
$ws.Range("I20").Value = 1700
$ws.Range("K20").Value = 1700
$ws.Range("M20").Value = -1470

# Row 35
$ws.Range("H35").Value = 1700
$ws.Range("I35").Value = 1700
$ws.Range("K35").Value = 1700
$ws.Range("M35").Value = -1321

# Row 86
$ws.Range("H86").Value = 4500
$ws.Range("I86").Value = 4500
$ws.Range("K86").Value = 4500
$ws.Range("M86").Value = -3377

# Row 89
$ws.Range("H89").Value = 4500
$ws.Range("I89").Value = 4500
$ws.Range("K89").Value = 22500
$ws.Range("M89").Value = -16884

# Row 141
$ws.Range("H141").Value = 5496.5
$ws.Range("J141").Value = 4996.3335
$ws.Range("L141").Value = 14989.0005
$ws.Range("N141").Value = -25349.0005


# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")

# Row 32
$ws.Range("I32").Value = 2060644.8
$ws.Range("K32").Value = 2060644.8
$ws.Range("M32").Value = -2060357.8

# Row 45
$ws.Range("H45").Value = 201438
$ws.Range("I45").Value = 201438
$ws.Range("K45").Value = 201438
$ws.Range("M45").Value = -201061

# Row 74
$ws.Range("H74").Value = 1325
$ws.Range("I74").Value = 1325
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1325
$ws.Range("L74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("N74").Value = -451

# Row 76
$ws.Range("H76").Value = 3899.6
$ws.Range("I76").Value = 5000
$ws.Range("J76").Value = 3624.5
$ws.Range("K76").Value = 5000
$ws.Range("L76").Value = 3624.5
$ws.Range("M76").Value = -4662
$ws.Range("N76").Value = -4300.5

# Row 77
$ws.Range("H77").Value = 1325
$ws.Range("I77").Value = 1325
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 6625
$ws.Range("L77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("N77").Value = -2257

# Row 79
$ws.Range("H79").Value = 3899.6
$ws.Range("I79").Value = 5000
$ws.Range("J79").Value = 3624.5
$ws.Range("K79").Value = 5000
$ws.Range("L79").Value = 3624.5
$ws.Range("M79").Value = -3830
$ws.Range("N79").Value = -5964.5

# Row 125
$ws.Range("H125").Value = 100000
$ws.Range("J125").Value = 100000
$ws.Range("L125").Value = 100000
$ws.Range("N125").Value = -109840


# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")

# Row 107
$ws.Range("H107").Value = 3301.1667
$ws.Range("I107").Value = 3301.1667
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 3301.1667
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -1381.1667


# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")

# Row 22
$ws.Range("H22").Value = 555.1429000000001
$ws.Range("I22").Value = 571.5
$ws.Range("J22").Value = 457
$ws.Range("K22").Value = 571.5
$ws.Range("L22").Value = 457
$ws.Range("M22").Value = -221.5
$ws.Range("N22").Value = -1157

# Row 31
$ws.Range("H31").Value = 1311.6666
$ws.Range("I31").Value = 1272.2858
$ws.Range("K31").Value = 1272.2858
$ws.Range("M31").Value = -977.2858000000001

# Row 34
$ws.Range("H34").Value = 1311.6666
$ws.Range("I34").Value = 1272.2858
$ws.Range("K34").Value = 1272.2858
$ws.Range("M34").Value = -1070.2858

# Row 94
$ws.Range("H94").Value = 163391.28
$ws.Range("I94").Value = 374336.66
$ws.Range("J94").Value = 5182.25
$ws.Range("K94").Value = 374336.66
$ws.Range("L94").Value = 5182.25
$ws.Range("M94").Value = -373885.66
$ws.Range("N94").Value = -6084.25

# Row 107
$ws.Range("H107").Value = 1362.7778
$ws.Range("I107").Value = 1329.6666
$ws.Range("J107").Value = 1429
$ws.Range("K107").Value = 1329.6666
$ws.Range("L107").Value = 1429
$ws.Range("M107").Value = 590.3334
$ws.Range("N107").Value = -5269

# Row 132
$ws.Range("H132").Value = 3307.3333
$ws.Range("I132").Value = 3529.6
$ws.Range("K132").Value = 10588.8
$ws.Range("M132").Value = -8058.799999999999


# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")

# Row 32
$ws.Range("H32").Value = 4995
$ws.Range("J32").Value = 4995
$ws.Range("L32").Value = 14985
$ws.Range("N32").Value = -15551

# Row 113
$ws.Range("H113").Value = 1557.6364
$ws.Range("I113").Value = 1613
$ws.Range("J113").Value = 1536.875
$ws.Range("K113").Value = 4839
$ws.Range("L113").Value = 4610.625
$ws.Range("M113").Value = -2669
$ws.Range("N113").Value = -8950.625

# Row 117
$ws.Range("H117").Value = 33509.668
$ws.Range("I117").Value = 264.5
$ws.Range("J117").Value = 100000
$ws.Range("K117").Value = 793.5
$ws.Range("L117").Value = 300000
$ws.Range("M117").Value = 2648.5
$ws.Range("N117").Value = -306884


# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")

# Row 80
$ws.Range("H80").Value = 3849.6667
$ws.Range("I80").Value = 4049
$ws.Range("J80").Value = 3750
$ws.Range("K80").Value = 4049
$ws.Range("L80").Value = 3750
$ws.Range("M80").Value = -3051
$ws.Range("N80").Value = -5746

# Row 83
$ws.Range("H83").Value = 3849.6667
$ws.Range("I83").Value = 4049
$ws.Range("J83").Value = 3750
$ws.Range("K83").Value = 20245
$ws.Range("L83").Value = 18750
$ws.Range("M83").Value = -15253
$ws.Range("N83").Value = -28734

# Row 97
$ws.Range("H97").Value = 2164.6667
$ws.Range("I97").Value = 747
$ws.Range("K97").Value = 747
$ws.Range("M97").Value = -251


# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")

# Row 22
$ws.Range("H22").Value = 1066
$ws.Range("J22").Value = 399
$ws.Range("L22").Value = 399
$ws.Range("N22").Value = -989

# Row 27
$ws.Range("H27").Value = 1066
$ws.Range("J27").Value = 399
$ws.Range("L27").Value = 399
$ws.Range("N27").Value = -613

# Row 68
$ws.Range("H68").Value = 2233.3333
$ws.Range("I68").Value = 2233.3333
$ws.Range("K68").Value = 2233.3333
$ws.Range("M68").Value = -1484.3333

# Row 71
$ws.Range("H71").Value = 2233.3333
$ws.Range("I71").Value = 2233.3333
$ws.Range("K71").Value = 11166.6665
$ws.Range("M71").Value = -7422.666499999999

# Row 132
$ws.Range("H132").Value = 6637.467
$ws.Range("I132").Value = 7854.8
$ws.Range("K132").Value = 23564.4
$ws.Range("M132").Value = -21034.4

# Row 136
$ws.Range("H136").Value = 2254.6155
$ws.Range("I136").Value = 1387.625
$ws.Range("J136").Value = 3641.8
$ws.Range("K136").Value = 4162.875
$ws.Range("L136").Value = 10925.4
$ws.Range("M136").Value = -1612.875
$ws.Range("N136").Value = -16025.4


# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")

# Row 41
$ws.Range("H41").Value = 19981.8
$ws.Range("I41").Value = 19978
$ws.Range("K41").Value = 19978
$ws.Range("M41").Value = -19588

# Row 75
$ws.Range("H75").Value = 74252
$ws.Range("J75").Value = 74000
$ws.Range("L75").Value = 74000
$ws.Range("N75").Value = -75872

# Row 78
$ws.Range("H78").Value = 74252
$ws.Range("J78").Value = 74000
$ws.Range("L78").Value = 222000
$ws.Range("N78").Value = -231360

# Row 132
$ws.Range("H132").Value = 1649.909
$ws.Range("I132").Value = 1620.4736
$ws.Range("J132").Value = 1836.3334
$ws.Range("K132").Value = 4861.4208
$ws.Range("L132").Value = 5509.0002
$ws.Range("M132").Value = -2331.4208
$ws.Range("N132").Value = -10569.0002

